$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix typo in the existing "Developmental Screening Gap Asessment" title
#    (this cell is currently row 13, column E) BEFORE inserting the new row,
#    so the shared-string table keeps the same append order as the target file.
$ws.Range("E13").Value2 = "Developmental Screening Gap Assessment - Prepared for First 5 Kern County"

# 2) Insert a new row at position 12 for the new May 2025 publication entry,
#    which pushes all the existing pro_report rows (old 12-22) down to 13-23.
$ws.Rows.Item(12).Insert()

# The inserted row copies formatting from the row above (row 11, a "conference"
# entry) which also populates stray P:S cells - clear those so only A:F remain.
$ws.Range("P12:S12").Clear()

# 3) Populate the new row 12 with the new publication entry.
$ws.Range("A12").Value2 = "pro_report"
$ws.Range("B12").Value2 = 1
$ws.Range("C12").Value2 = "Wright, J."
$ws.Range("D12").Value2 = 2025
$ws.Range("E12").Value2 = "Parent/Guardian Survey on Childhood Immunization and Vaccination: Spring 2025 - Prepared for the Kern County Immunization Coalition"
$ws.Range("F12").Value2 = "Applied Survey Research"
$ws.Rows.Item(12).RowHeight = 119

# 4) Renumber the "order" column (B) for the rows that shifted down one
#    position (old order 1-11 in rows 13-23 becomes 2-12).
$ws.Range("B13").Value2 = 2
$ws.Range("B14").Value2 = 3
$ws.Range("B15").Value2 = 4
$ws.Range("B16").Value2 = 5
$ws.Range("B17").Value2 = 6
$ws.Range("B18").Value2 = 7
$ws.Range("B19").Value2 = 8
$ws.Range("B20").Value2 = 9
$ws.Range("B21").Value2 = 10
$ws.Range("B22").Value2 = 11
$ws.Range("B23").Value2 = 12

# 5) Update the sheet view: select column K (scrolling back to the top-left)
#    instead of the previous scrolled-down B30 selection.
[void]$ws.Columns.Item(11).Select()
